# Apply the rating / rate_updated_at corrections described by the diff for
# the ratings_only query-results export (rows 19-22 of the sheet).
#
# Row 19: rating 5 -> 3           (rate_updated_at text unchanged cell-wise,
#                                   underlying timestamp text updated below)
# Row 20: rating 4 -> 5
# Row 21: rating 4 -> 7,          rate_updated_at -> 2024-12-12 13:49:33
# Row 22: rating "10" -> "7"      (stored as text, same as before)
#
# rate_updated_at (column K) values for rows 19 & 20 are refreshed in place;
# row 21's timestamp becomes the same later value already used by row 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column J : rating ------------------------------------------------
$ws.Range("J19").Value = 3
$ws.Range("J20").Value = 5
$ws.Range("J21").Value = 7

# J22 historically stores its rating as text ("10" -> "7"), not a number,
# so force a Text format before writing it to keep it a text value instead
# of Excel auto-converting it to the number 7.
$ws.Range("J22").NumberFormat = "@"
$ws.Range("J22").Value = "7"

# --- column K : rate_updated_at ---------------------------------------
$ws.Range("K19").Value = "2024-12-12 10:25:26"
$ws.Range("K20").Value = "2024-12-12 13:47:28"
$ws.Range("K21").Value = "2024-12-12 13:49:33"

# Row 22 used to share its rate_updated_at string with row 21
# ("2024-12-12 10:10:33"); that shared text is renamed to
# "2024-12-12 13:49:33" by the edit, so row 22's displayed value moves too.
$ws.Range("K22").Value = "2024-12-12 13:49:33"
